# Append the new survey response row (row 28) coming from the
# "xiaoyi" Streamlit SmartScore submission.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28

# --- plain text / identifier columns -------------------------------------
$ws.Cells.Item($row, 1).Value  = "xiaoyi_20251202_134614"      # A: ID_Participante

# B: Grupo_Experimental -> empty string (quote-prefixed empty text, not a blank cell)
$ws.Cells.Item($row, 2).Value  = "'"
$ws.Cells.Item($row, 2).Style  = "Normal"

$ws.Cells.Item($row, 3).Value  = "xiaoyi"                       # C: Nombre Completo
$ws.Cells.Item($row, 4).Value  = 26                              # D: Edad (number)
$ws.Cells.Item($row, 5).Value  = "Female"                       # E: Género
$ws.Cells.Item($row, 6).Value  = "2025-12-02 13:46:14"          # F: Fecha (literal text)

# G: Pesos -> multi-line JSON blob, stored verbatim as text
$pesos = "{" + [char]10 +
  "  ""portion"": 0.4," + [char]10 +
  "  ""diet"": 0.7142857142857143," + [char]10 +
  "  ""salt"": 0.2," + [char]10 +
  "  ""fat"": 0.8," + [char]10 +
  "  ""natural"": 0.2," + [char]10 +
  "  ""convenience"": 1.0," + [char]10 +
  "  ""price"": 1.0" + [char]10 +
  "}"
$ws.Cells.Item($row, 7).Value = $pesos

# --- Instant Noodles ---------------------------------------------------
$ws.Cells.Item($row, 8).Value   = "Maruchan Ramen Sabor Pollo"
$ws.Cells.Item($row, 9).Value   = "'0.591"
$ws.Cells.Item($row, 9).Style   = "Normal"
$ws.Cells.Item($row, 10).Value  = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"

$ws.Cells.Item($row, 11).Value  = "Nongshim Neoguri Spicy Seafood"
$ws.Cells.Item($row, 12).Value  = "'0.532"
$ws.Cells.Item($row, 12).Style  = "Normal"
$ws.Cells.Item($row, 13).Value  = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"

$ws.Cells.Item($row, 14).Value  = "Nissin Chow Mein Teriyaki Beef"
$ws.Cells.Item($row, 15).Value  = "'0.518"
$ws.Cells.Item($row, 15).Style  = "Normal"
$ws.Cells.Item($row, 16).Value  = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"

# --- Mac & Cheese --------------------------------------------------------
$ws.Cells.Item($row, 17).Value  = "Velveeta Original Shells & Cheese (microwave cups)"
$ws.Cells.Item($row, 18).Value  = "'0.607"
$ws.Cells.Item($row, 18).Style  = "Normal"
$ws.Cells.Item($row, 19).Value  = "Muy cremoso, porción individual, rápido, salado, ideal para niños"

$ws.Cells.Item($row, 20).Value  = "Kraft Macaroni & Cheese Dinner"
$ws.Cells.Item($row, 21).Value  = "'0.520"
$ws.Cells.Item($row, 21).Style  = "Normal"
$ws.Cells.Item($row, 22).Value  = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"

$ws.Cells.Item($row, 23).Value  = "Amy’s Macaroni & Cheese (frozen)"
$ws.Cells.Item($row, 24).Value  = "'0.443"
$ws.Cells.Item($row, 24).Style  = "Normal"
$ws.Cells.Item($row, 25).Value  = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"

# --- Ready to Eat ----------------------------------------------------------
$ws.Cells.Item($row, 26).Value  = "StarKist Chicken Creations (Chicken Salad)"
$ws.Cells.Item($row, 27).Value  = "'0.718"
$ws.Cells.Item($row, 27).Style  = "Normal"
$ws.Cells.Item($row, 28).Value  = "Portátil, saludable, fácil, buena textura, sabor suave"

$ws.Cells.Item($row, 29).Value  = "Jack Link’s Beef Jerky Original"
$ws.Cells.Item($row, 30).Value  = "'0.705"
$ws.Cells.Item($row, 30).Style  = "Normal"
$ws.Cells.Item($row, 31).Value  = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

$ws.Cells.Item($row, 32).Value  = "Wild Planet Wild Tuna Pasta Salad"
$ws.Cells.Item($row, 33).Value  = "'0.673"
$ws.Cells.Item($row, 33).Style  = "Normal"
$ws.Cells.Item($row, 34).Value  = "Sabor fresco, buena proteína, saludable, porción algo pequeña"

# Row height stays default — undo the autofit bump triggered by the
# multi-line JSON cell above so the row matches the rest of the sheet.
$ws.Rows.Item($row).AutoFit()
